# Update the "final data" values for rows 2-6 (Bridgetown, L&D, Levy-shipping,
# NCQG, Levy-aviation) in columns B (All), K (Japan), L (Russia) and N (USA).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.793366965036613
$ws.Range("K2").Value = 0.805037944194611
$ws.Range("L2").Value = 0.829224529123518
$ws.Range("N2").Value = 0.740197166164522

$ws.Range("B3").Value = 0.748295871658231
$ws.Range("K3").Value = 0.727293734869913
$ws.Range("L3").Value = 0.866768469161838
$ws.Range("N3").Value = 0.697898274519676

$ws.Range("B4").Value = 0.699644803977167
$ws.Range("K4").Value = 0.58767644567848
$ws.Range("L4").Value = 0.72514029245591
$ws.Range("N4").Value = 0.674357186299684

$ws.Range("B5").Value = 0.683174203642518
$ws.Range("K5").Value = 0.59008286280892
$ws.Range("L5").Value = 0.875589503017251
$ws.Range("N5").Value = 0.614162213651268

$ws.Range("B6").Value = 0.526319576644156
$ws.Range("K6").Value = 0.46350801136536
$ws.Range("L6").Value = 0.514246894490421
$ws.Range("N6").Value = 0.507370211647722
